$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 412.53333
$ws.Range("I39").Value = 318.5
$ws.Range("J39").Value = 600.6
$ws.Range("K39").Value = 955.5
$ws.Range("L39").Value = 1801.8
$ws.Range("M39").Value = -659.5
$ws.Range("N39").Value = -2393.8

$ws.Range("H64").Value = 125004900
$ws.Range("I64").Value = 166672060
$ws.Range("J64").Value = 3400
$ws.Range("K64").Value = 166672060
$ws.Range("L64").Value = 3400
$ws.Range("M64").Value = -166671812
$ws.Range("N64").Value = -3896

$ws.Range("H67").Value = 125004900
$ws.Range("I67").Value = 166672060
$ws.Range("J67").Value = 3400
$ws.Range("K67").Value = 166672060
$ws.Range("L67").Value = 3400
$ws.Range("M67").Value = -166671202
$ws.Range("N67").Value = -5116

$ws.Range("H100").Value = 3158.348
$ws.Range("I100").Value = 2932.7856
$ws.Range("J100").Value = 3509.2222
$ws.Range("K100").Value = 2932.7856
$ws.Range("L100").Value = 3509.2222
$ws.Range("M100").Value = -2391.7856
$ws.Range("N100").Value = -4591.2222

$ws.Range("H103").Value = 585
$ws.Range("I103").Value = 998
$ws.Range("J103").Value = 447.33334
$ws.Range("K103").Value = 2994
$ws.Range("L103").Value = 1342.00002
$ws.Range("M103").Value = -2408
$ws.Range("N103").Value = -2514.00002

$ws.Range("H106").Value = 5476.6665
$ws.Range("I106").Value = 6127.5713
$ws.Range("J106").Value = 3198.5
$ws.Range("K106").Value = 6127.5713
$ws.Range("L106").Value = 3198.5
$ws.Range("M106").Value = -5496.5713
$ws.Range("N106").Value = -4460.5

$ws.Range("H116").Value = 7638.5
$ws.Range("I116").Value = 7340.7144
$ws.Range("J116").Value = 8333.333000000001
$ws.Range("K116").Value = 7340.7144
$ws.Range("L116").Value = 8333.333000000001
$ws.Range("M116").Value = -3898.7144
$ws.Range("N116").Value = -15217.333

$ws.Range("H138").Value = 2328.4148
$ws.Range("I138").Value = 1774.6666
$ws.Range("J138").Value = 2518.2715
$ws.Range("K138").Value = 5323.9998
$ws.Range("L138").Value = 7554.814499999999
$ws.Range("M138").Value = -183.9997999999996
$ws.Range("N138").Value = -17834.8145

$ws.Range("H141").Value = 3826.7222
$ws.Range("I141").Value = 4113.0713
$ws.Range("J141").Value = 2824.5
$ws.Range("K141").Value = 12339.2139
$ws.Range("L141").Value = 8473.5
$ws.Range("M141").Value = -7159.213899999999
$ws.Range("N141").Value = -18833.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4499
$ws.Range("I2").Value = 4499
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 4499
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -1506.8

$ws.Range("H45").Value = 2920.2856
$ws.Range("I45").Value = 2793.524
$ws.Range("J45").Value = 3300.5715
$ws.Range("K45").Value = 2793.524
$ws.Range("L45").Value = 3300.5715
$ws.Range("M45").Value = -2416.524
$ws.Range("N45").Value = -4054.5715

$ws.Range("H102").Value = 8168.8
$ws.Range("I102").Value = 3616.3333
$ws.Range("J102").Value = 14997.5
$ws.Range("K102").Value = 3616.3333
$ws.Range("L102").Value = 14997.5
$ws.Range("M102").Value = -1994.3333
$ws.Range("N102").Value = -18241.5

$ws.Range("H116").Value = 4499
$ws.Range("I116").Value = 4499
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 4499
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 674.2

$ws.Range("H122").Value = 1538550.5
$ws.Range("I122").Value = 2440221.2
$ws.Range("J122").Value = 5710.4
$ws.Range("K122").Value = 7320663.600000001
$ws.Range("L122").Value = 17131.2
$ws.Range("M122").Value = -7318213.600000001
$ws.Range("N122").Value = -22031.2

$ws.Range("H132").Value = 2510929.2
$ws.Range("I132").Value = 3444.3845
$ws.Range("J132").Value = 7167686.5
$ws.Range("K132").Value = 10333.1535
$ws.Range("L132").Value = 21503059.5
$ws.Range("M132").Value = -7803.1535
$ws.Range("N132").Value = -21508119.5

$ws.Range("N2").ClearContents()
$ws.Range("N116").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4499
$ws.Range("I3").Value = 4499
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 4499
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1505.8

$ws.Range("H20").Value = 6076400
$ws.Range("I20").Value = 11115332
$ws.Range("J20").Value = 29681.92
$ws.Range("K20").Value = 11115332
$ws.Range("L20").Value = 29681.92
$ws.Range("M20").Value = -11115085
$ws.Range("N20").Value = -30175.92

$ws.Range("H99").Value = 15291.073
$ws.Range("I99").Value = 16264.968
$ws.Range("J99").Value = 12272
$ws.Range("K99").Value = 16264.968
$ws.Range("L99").Value = 12272
$ws.Range("M99").Value = -14766.968
$ws.Range("N99").Value = -15268

$ws.Range("H105").Value = 111112730
$ws.Range("I105").Value = 111112730
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 111112730
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -111110983

$ws.Range("N3").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 11100.2
$ws.Range("I58").Value = 4075.0293
$ws.Range("J58").Value = 26028.688
$ws.Range("K58").Value = 4075.0293
$ws.Range("L58").Value = 26028.688
$ws.Range("M58").Value = -3872.0293
$ws.Range("N58").Value = -26434.688

$ws.Range("H62").Value = 2818.75
$ws.Range("I62").Value = 2758.3333
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 2758.3333
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -2134.3333
$ws.Range("N62").Value = -4248

$ws.Range("H65").Value = 2818.75
$ws.Range("I65").Value = 2758.3333
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 13791.6665
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -10671.6665
$ws.Range("N65").Value = -21240

$ws.Range("H105").Value = 8662.214
$ws.Range("I105").Value = 9226.23
$ws.Range("J105").Value = 1330
$ws.Range("K105").Value = 9226.23
$ws.Range("L105").Value = 1330
$ws.Range("M105").Value = -7479.23
$ws.Range("N105").Value = -4824

$ws.Range("H132").Value = 2589.2
$ws.Range("I132").Value = 2202.7144
$ws.Range("J132").Value = 8000
$ws.Range("K132").Value = 6608.1432
$ws.Range("L132").Value = 24000
$ws.Range("M132").Value = -4078.1432
$ws.Range("N132").Value = -29060

$ws.Range("H134").Value = 25646160
$ws.Range("I134").Value = 1644.9584
$ws.Range("J134").Value = 66677384
$ws.Range("K134").Value = 4934.8752
$ws.Range("L134").Value = 200032152
$ws.Range("M134").Value = -2399.8752
$ws.Range("N134").Value = -200037222

$ws.Range("H136").Value = 11100.2
$ws.Range("I136").Value = 4075.0293
$ws.Range("J136").Value = 26028.688
$ws.Range("K136").Value = 12225.0879
$ws.Range("L136").Value = 78086.064
$ws.Range("M136").Value = -9675.0879
$ws.Range("N136").Value = -83186.064

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 53.52
$ws.Range("I2").Value = 51.636364
$ws.Range("J2").Value = 55
$ws.Range("K2").Value = 309.818184
$ws.Range("L2").Value = 330
$ws.Range("M2").Value = -196.818184
$ws.Range("N2").Value = -556

$ws.Range("H131").Value = 1454.4
$ws.Range("I131").Value = 1030
$ws.Range("J131").Value = 1458.6869
$ws.Range("K131").Value = 3090
$ws.Range("L131").Value = 4376.0607
$ws.Range("M131").Value = 1950
$ws.Range("N131").Value = -14456.0607

$ws.Range("H132").Value = 1888.1818
$ws.Range("I132").Value = 1826.75
$ws.Range("J132").Value = 1961.9
$ws.Range("K132").Value = 16440.75
$ws.Range("L132").Value = 17657.1
$ws.Range("M132").Value = -13910.75
$ws.Range("N132").Value = -22717.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 23888.6
$ws.Range("I80").Value = 27972
$ws.Range("J80").Value = 21166.334
$ws.Range("K80").Value = 27972
$ws.Range("L80").Value = 21166.334
$ws.Range("M80").Value = -26974
$ws.Range("N80").Value = -23162.334

$ws.Range("H83").Value = 23888.6
$ws.Range("I83").Value = 27972
$ws.Range("J83").Value = 21166.334
$ws.Range("K83").Value = 139860
$ws.Range("L83").Value = 105831.67
$ws.Range("M83").Value = -134868
$ws.Range("N83").Value = -115815.67

$ws.Range("H113").Value = 3962.6
$ws.Range("I113").Value = 3285
$ws.Range("J113").Value = 4209
$ws.Range("K113").Value = 3285
$ws.Range("L113").Value = 4209
$ws.Range("M113").Value = -1115
$ws.Range("N113").Value = -8549

$ws.Range("H132").Value = 3545.625
$ws.Range("I132").Value = 1430.7273
$ws.Range("J132").Value = 8198.4
$ws.Range("K132").Value = 4292.1819
$ws.Range("L132").Value = 24595.2
$ws.Range("M132").Value = -1762.1819
$ws.Range("N132").Value = -29655.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1994792.2
$ws.Range("I7").Value = 3675797.8
$ws.Range("J7").Value = 8149.364
$ws.Range("K7").Value = 3675797.8
$ws.Range("L7").Value = 8149.364
$ws.Range("M7").Value = -3675685.8
$ws.Range("N7").Value = -8373.364

$ws.Range("H22").Value = 25002002
$ws.Range("I22").Value = 1540.3334
$ws.Range("J22").Value = 52634092
$ws.Range("K22").Value = 1540.3334
$ws.Range("L22").Value = 52634092
$ws.Range("M22").Value = -1245.3334
$ws.Range("N22").Value = -52634682

$ws.Range("H27").Value = 25002002
$ws.Range("I27").Value = 1540.3334
$ws.Range("J27").Value = 52634092
$ws.Range("K27").Value = 1540.3334
$ws.Range("L27").Value = 52634092
$ws.Range("M27").Value = -1433.3334
$ws.Range("N27").Value = -52634306

$ws.Range("H40").Value = 3684500.8
$ws.Range("I40").Value = 9060.625
$ws.Range("J40").Value = 7359941
$ws.Range("K40").Value = 9060.625
$ws.Range("L40").Value = 7359941
$ws.Range("M40").Value = -8924.625
$ws.Range("N40").Value = -7360213

$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 932.5

$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 4662.5

$ws.Range("H82").Value = 2910.8635
$ws.Range("I82").Value = 2945.0588
$ws.Range("J82").Value = 2794.6
$ws.Range("K82").Value = 2945.0588
$ws.Range("L82").Value = 2794.6
$ws.Range("M82").Value = -2584.0588
$ws.Range("N82").Value = -3516.6

$ws.Range("H85").Value = 2910.8635
$ws.Range("I85").Value = 2945.0588
$ws.Range("J85").Value = 2794.6
$ws.Range("K85").Value = 2945.0588
$ws.Range("L85").Value = 2794.6
$ws.Range("M85").Value = -1697.0588
$ws.Range("N85").Value = -5290.6

$ws.Range("H126").Value = 1994792.2
$ws.Range("I126").Value = 3675797.8
$ws.Range("J126").Value = 8149.364
$ws.Range("K126").Value = 11027393.4
$ws.Range("L126").Value = 24448.092
$ws.Range("M126").Value = -11024923.4
$ws.Range("N126").Value = -29388.092

$ws.Range("H136").Value = 17231.654
$ws.Range("I136").Value = 14751.125
$ws.Range("J136").Value = 21200.5
$ws.Range("K136").Value = 44253.375
$ws.Range("L136").Value = 63601.5
$ws.Range("M136").Value = -41703.375
$ws.Range("N136").Value = -68701.5

$ws.Range("M68").ClearContents()
$ws.Range("N68").ClearContents()
$ws.Range("M71").ClearContents()
$ws.Range("N71").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 522.3913
$ws.Range("I100").Value = 489.6154
$ws.Range("J100").Value = 565
$ws.Range("K100").Value = 979.2308
$ws.Range("L100").Value = 1130
$ws.Range("M100").Value = -438.2308
$ws.Range("N100").Value = -2212

$ws.Range("H132").Value = 11491.942
$ws.Range("I132").Value = 3711.4883
$ws.Range("J132").Value = 48665.223
$ws.Range("K132").Value = 11134.4649
$ws.Range("L132").Value = 145995.669
$ws.Range("M132").Value = -8604.464899999999
$ws.Range("N132").Value = -151055.669

$ws.Range("H136").Value = 11896.426
$ws.Range("I136").Value = 3159.2122
$ws.Range("J136").Value = 32491.285
$ws.Range("K136").Value = 9477.6366
$ws.Range("L136").Value = 97473.855
$ws.Range("M136").Value = -6927.6366
$ws.Range("N136").Value = -102573.855
